$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A23: 2024-02-06 -> 2024-02-07 -------------------------------------
# A plain ".Value = '2024-02-07'" assignment gets auto-detected as a date
# literal and stored as a numeric serial (this column is meant to stay
# plain text, matching every other row in the sheet). Route the literal
# through a text formula and flatten it back to a static value via
# copy/paste-special so the cell ends up as a normal string cell (no
# lingering formula, no extra number-format/style picked up along the way).
$ws.Range("A23").Formula = '="2024-02-07"'
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = $false

# --- B23: "Created: Feat: Test PR" -> multi-line note -------------------
$ws.Range("B23").Value = "Created: Fix: Action filters default enabled, Fix: All actions ticked by default`nClosed: Fix: All actions ticked by default"

# Re-fit the row so the newly multi-line B23 doesn't leave an explicit
# custom row height behind (row 22, which already holds multi-line text,
# has none either).
$ws.Rows(23).AutoFit() | Out-Null
